$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Add a "submit" row so the form can be submitted.
$ws.Range("C12").Value = "submit"
$ws.Range("D12").Value = "bla"
$ws.Range("F12").Value = "Go on!"

# The "instruction" row is now a "note" row - notes are always displayed
# (treated like answering any other item), so its XLSForm type changes
# from "instruction" to "note".
$ws.Range("C2").Value = "note"

# Add a trailing "note" row shown after submission.
$ws.Range("F13").Value = "Good work, chap!"
$ws.Range("D13").Value = "xx"
$ws.Range("C13").Value = "note"

# Match styling used by the rest of the data rows (wrap text).
$ws.Range("C12").WrapText = $true
$ws.Range("D12").WrapText = $true
$ws.Range("F12").WrapText = $true
$ws.Range("C13").WrapText = $true
$ws.Range("D13").WrapText = $true
$ws.Range("F13").WrapText = $true

# Update the stored cursor/selection position on the survey sheet.
$ws.Range("C3").Select()
